$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M15").Value = -2783140.4
$ws.Range("K15").Value = 2783309.4
$ws.Range("H15").Value = 927769.8
$ws.Range("I15").Value = 927769.8
$ws.Range("H32").Value = 1520.375
$ws.Range("L32").Value = 1451.8572
$ws.Range("J32").Value = 1451.8572
$ws.Range("N32").Value = -2103.8572
$ws.Range("L38").Value = 15792.429
$ws.Range("H38").Value = 1918.8
$ws.Range("J38").Value = 5264.143
$ws.Range("N38").Value = -16536.429
$ws.Range("J40").Value = 4580
$ws.Range("N40").Value = -4930
$ws.Range("L40").Value = 4580
$ws.Range("H40").Value = 4008.1304
$ws.Range("I43").Value = 2800.3333
$ws.Range("J43").Value = 4811.75
$ws.Range("N43").Value = -4949.75
$ws.Range("M43").Value = -2731.3333
$ws.Range("K43").Value = 2800.3333
$ws.Range("L43").Value = 4811.75
$ws.Range("H43").Value = 4263.1816
$ws.Range("K55").Value = 116.8
$ws.Range("L55").Value = 67.5
$ws.Range("H55").Value = 102.71429
$ws.Range("I55").Value = 116.8
$ws.Range("J55").Value = 67.5
$ws.Range("N55").Value = -495.5
$ws.Range("M55").Value = 97.2
$ws.Range("H62").Value = 6481
$ws.Range("I62").Value = 2499.5
$ws.Range("M62").Value = -1875.5
$ws.Range("K62").Value = 2499.5
$ws.Range("H65").Value = 6481
$ws.Range("I65").Value = 2499.5
$ws.Range("M65").Value = -9377.5
$ws.Range("K65").Value = 12497.5
$ws.Range("J69").Value = 19799.8
$ws.Range("N69").Value = -61147.39999999999
$ws.Range("L69").Value = 59399.39999999999
$ws.Range("H69").Value = 17999.723
$ws.Range("J72").Value = 19799.8
$ws.Range("N72").Value = -186934.2
$ws.Range("L72").Value = 178198.2
$ws.Range("H72").Value = 17999.723
$ws.Range("K80").Value = 2353.5
$ws.Range("H80").Value = 1032.3334
$ws.Range("I80").Value = 784.5
$ws.Range("M80").Value = -1355.5
$ws.Range("K83").Value = 7060.5
$ws.Range("H83").Value = 1032.3334
$ws.Range("I83").Value = 784.5
$ws.Range("M83").Value = -2068.5
$ws.Range("K106").Value = 44001572
$ws.Range("H106").Value = 33848264
$ws.Range("I106").Value = 44001572
$ws.Range("M106").Value = -44000941
$ws.Range("L112").Value = 13097.5005
$ws.Range("H112").Value = 3799.25
$ws.Range("J112").Value = 4365.8335
$ws.Range("N112").Value = -15313.5005
$ws.Range("H132").Value = 2889.3704
$ws.Range("I132").Value = 2683.9363
$ws.Range("M132").Value = -5521.8089
$ws.Range("K132").Value = 8051.8089
$ws.Range("K135").Value = 25191
$ws.Range("H135").Value = 8655.333000000001
$ws.Range("I135").Value = 2799
$ws.Range("M135").Value = -22656
$ws.Range("M137").Value = -10856.0772
$ws.Range("J137").Value = 8160.684
$ws.Range("N137").Value = -29582.052
$ws.Range("K137").Value = 13406.0772
$ws.Range("L137").Value = 24482.052
$ws.Range("H137").Value = 6027.533
$ws.Range("I137").Value = 4468.6924
$ws.Range("K138").Value = 13429.125
$ws.Range("L138").Value = 19882.5
$ws.Range("H138").Value = 5737.3794
$ws.Range("I138").Value = 4476.375
$ws.Range("J138").Value = 6627.5
$ws.Range("N138").Value = -30162.5
$ws.Range("M138").Value = -8289.125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M2").Value = -1254.875
$ws.Range("K2").Value = 1367.875
$ws.Range("H2").Value = 1290.4546
$ws.Range("I2").Value = 1367.875
$ws.Range("M5").Value = -267.7143
$ws.Range("K5").Value = 379.7143
$ws.Range("H5").Value = 1000.2727
$ws.Range("I5").Value = 379.7143
$ws.Range("H32").Value = 5616.0215
$ws.Range("K32").Value = 3749.068
$ws.Range("I32").Value = 3749.068
$ws.Range("M32").Value = -3462.068
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 10000
$ws.Range("H41").Value = 7500
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 10000
$ws.Range("N41").Value = -10828
$ws.Range("M41").Value = -4586
$ws.Range("N44").ClearContents()
$ws.Range("L44").Value = 0
$ws.Range("H44").Value = 5000
$ws.Range("J44").Value = 0
$ws.Range("K45").Value = 3999
$ws.Range("H45").Value = 4599.2
$ws.Range("I45").Value = 3999
$ws.Range("M45").Value = -3622
$ws.Range("L55").Value = 57749.668
$ws.Range("H55").Value = 57749.668
$ws.Range("J55").Value = 57749.668
$ws.Range("N55").Value = -58379.668
$ws.Range("K61").Value = 3321.375
$ws.Range("H61").Value = 29418092
$ws.Range("I61").Value = 3321.375
$ws.Range("M61").Value = -3109.375
$ws.Range("J63").Value = 4383.1665
$ws.Range("N63").Value = -5755.1665
$ws.Range("L63").Value = 4383.1665
$ws.Range("H63").Value = 4328.4287
$ws.Range("J66").Value = 4383.1665
$ws.Range("N66").Value = -28779.8325
$ws.Range("L66").Value = 21915.8325
$ws.Range("H66").Value = 4328.4287
$ws.Range("K74").Value = 2643.8333
$ws.Range("H74").Value = 4861.773
$ws.Range("I74").Value = 2643.8333
$ws.Range("M74").Value = -1769.8333
$ws.Range("K77").Value = 13219.1665
$ws.Range("H77").Value = 4861.773
$ws.Range("I77").Value = 2643.8333
$ws.Range("M77").Value = -8851.166499999999
$ws.Range("H97").Value = 748.2857
$ws.Range("I97").Value = 728.9231
$ws.Range("M97").Value = -232.9231
$ws.Range("K97").Value = 728.9231
$ws.Range("J104").Value = 158684.38
$ws.Range("N104").Value = -165672.38
$ws.Range("L104").Value = 158684.38
$ws.Range("H104").Value = 158684.38
$ws.Range("H116").Value = 1290.4546
$ws.Range("I116").Value = 1367.875
$ws.Range("M116").Value = 926.125
$ws.Range("K116").Value = 1367.875
$ws.Range("H136").Value = 29418092
$ws.Range("I136").Value = 3321.375
$ws.Range("M136").Value = -7414.125
$ws.Range("K136").Value = 9964.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K3").Value = 1367.875
$ws.Range("H3").Value = 1290.4546
$ws.Range("I3").Value = 1367.875
$ws.Range("M3").Value = -1253.875
$ws.Range("I4").Value = 379.7143
$ws.Range("M4").Value = -264.7143
$ws.Range("K4").Value = 379.7143
$ws.Range("H4").Value = 1000.2727
$ws.Range("H26").Value = 39997
$ws.Range("I26").Value = 39997
$ws.Range("M26").Value = -39705
$ws.Range("K26").Value = 39997
$ws.Range("N28").ClearContents()
$ws.Range("L28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L35").Value = 63937.5
$ws.Range("H35").Value = 63937.5
$ws.Range("J35").Value = 63937.5
$ws.Range("N35").Value = -64557.5
$ws.Range("K80").Value = 75353
$ws.Range("L80").Value = 14568.667
$ws.Range("H80").Value = 33271.54
$ws.Range("I80").Value = 75353
$ws.Range("M80").Value = -74355
$ws.Range("J80").Value = 14568.667
$ws.Range("N80").Value = -16564.667
$ws.Range("K83").Value = 376765
$ws.Range("L83").Value = 72843.33499999999
$ws.Range("H83").Value = 33271.54
$ws.Range("I83").Value = 75353
$ws.Range("J83").Value = 14568.667
$ws.Range("N83").Value = -82827.33499999999
$ws.Range("M83").Value = -371773
$ws.Range("H107").Value = 862.8570999999999
$ws.Range("J107").Value = 2700
$ws.Range("N107").Value = -6540
$ws.Range("L107").Value = 2700

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M31").Value = -2590.9768
$ws.Range("K31").Value = 2885.9768
$ws.Range("H31").Value = 4927.492
$ws.Range("I31").Value = 2885.9768
$ws.Range("M34").Value = -2683.9768
$ws.Range("K34").Value = 2885.9768
$ws.Range("H34").Value = 4927.492
$ws.Range("I34").Value = 2885.9768
$ws.Range("K58").Value = 3236.923
$ws.Range("H58").Value = 6822.35
$ws.Range("L58").Value = 13481
$ws.Range("I58").Value = 3236.923
$ws.Range("J58").Value = 13481
$ws.Range("N58").Value = -13887
$ws.Range("M58").Value = -3033.923
$ws.Range("H68").Value = 158156.12
$ws.Range("I68").Value = 88000
$ws.Range("M68").Value = -87251
$ws.Range("J68").Value = 168178.42
$ws.Range("N68").Value = -169676.42
$ws.Range("K68").Value = 88000
$ws.Range("L68").Value = 168178.42
$ws.Range("L71").Value = 504535.26
$ws.Range("H71").Value = 158156.12
$ws.Range("I71").Value = 88000
$ws.Range("J71").Value = 168178.42
$ws.Range("M71").Value = -260256
$ws.Range("N71").Value = -512023.26
$ws.Range("K71").Value = 264000
$ws.Range("H122").Value = 1328.2
$ws.Range("K122").Value = 4020.6921
$ws.Range("L122").Value = 3750
$ws.Range("I122").Value = 1340.2307
$ws.Range("M122").Value = -1570.6921
$ws.Range("J122").Value = 1250
$ws.Range("N122").Value = -8650
$ws.Range("H132").Value = 5254.385
$ws.Range("I132").Value = 4118.8184
$ws.Range("M132").Value = -9826.4552
$ws.Range("J132").Value = 11500
$ws.Range("N132").Value = -39560
$ws.Range("K132").Value = 12356.4552
$ws.Range("L132").Value = 34500
$ws.Range("H136").Value = 6822.35
$ws.Range("I136").Value = 3236.923
$ws.Range("J136").Value = 13481
$ws.Range("N136").Value = -45543
$ws.Range("M136").Value = -7160.769
$ws.Range("K136").Value = 9710.769
$ws.Range("L136").Value = 40443
$ws.Range("L141").Value = 157999.4
$ws.Range("H141").Value = 139999.5
$ws.Range("J141").Value = 157999.4
$ws.Range("N141").Value = -168359.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M2").Value = -618.000016
$ws.Range("K2").Value = 731.000016
$ws.Range("H2").Value = 123.84615
$ws.Range("I2").Value = 121.833336
$ws.Range("L39").Value = 24316.3329
$ws.Range("H39").Value = 8105.4443
$ws.Range("J39").Value = 8105.4443
$ws.Range("N39").Value = -24904.3329
$ws.Range("H87").Value = 750
$ws.Range("I87").Value = 750
$ws.Range("M87").Value = -1002
$ws.Range("K87").Value = 2250
$ws.Range("K90").Value = 6750
$ws.Range("H90").Value = 750
$ws.Range("I90").Value = 750
$ws.Range("M90").Value = -510
$ws.Range("L96").Value = 5097
$ws.Range("H96").Value = 2065.7778
$ws.Range("J96").Value = 1699
$ws.Range("N96").Value = -9215
$ws.Range("K131").Value = 2096.25
$ws.Range("L131").Value = 6823917
$ws.Range("H131").Value = 758678.9
$ws.Range("I131").Value = 698.75
$ws.Range("J131").Value = 2274639
$ws.Range("N131").Value = -6833997
$ws.Range("M131").Value = 2943.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I11").Value = 1285499.9
$ws.Range("M11").Value = -1285360.9
$ws.Range("K11").Value = 1285499.9
$ws.Range("H11").Value = 1169540.8
$ws.Range("I14").Value = 212941.5
$ws.Range("J14").Value = 9119.799999999999
$ws.Range("M14").Value = -212773.5
$ws.Range("N14").Value = -9455.799999999999
$ws.Range("K14").Value = 212941.5
$ws.Range("L14").Value = 9119.799999999999
$ws.Range("H14").Value = 152993.94
$ws.Range("L19").Value = 20000
$ws.Range("H19").Value = 12500
$ws.Range("J19").Value = 20000
$ws.Range("N19").Value = -20576
$ws.Range("H20").Value = 39989
$ws.Range("I20").Value = 39989
$ws.Range("J20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("M20").Value = -39744
$ws.Range("K20").Value = 39989
$ws.Range("L20").Value = 0
$ws.Range("L54").Value = 24975
$ws.Range("H54").Value = 24975
$ws.Range("J54").Value = 24975
$ws.Range("N54").Value = -25755
$ws.Range("H107").Value = 458.8
$ws.Range("I107").Value = 348.5
$ws.Range("J107").Value = 900
$ws.Range("N107").Value = -4740
$ws.Range("M107").Value = 1571.5
$ws.Range("K107").Value = 348.5
$ws.Range("L107").Value = 900

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 2214.2856
$ws.Range("L6").Value = 80027.336
$ws.Range("H6").Value = 80027.336
$ws.Range("J6").Value = 80027.336
$ws.Range("N6").Value = -80251.336
$ws.Range("H15").Value = 2214.2856
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("H46").Value = 5632.1665
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 5948.25
$ws.Range("N46").Value = -6324.25
$ws.Range("M46").Value = -4812
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 5948.25
$ws.Range("K58").Value = 2283.3333
$ws.Range("H58").Value = 11719.8
$ws.Range("L58").Value = 25874.5
$ws.Range("I58").Value = 2283.3333
$ws.Range("J58").Value = 25874.5
$ws.Range("N58").Value = -26394.5
$ws.Range("M58").Value = -2023.3333
$ws.Range("H68").Value = 6138.8
$ws.Range("I68").Value = 6173.5
$ws.Range("M68").Value = -5424.5
$ws.Range("K68").Value = 6173.5
$ws.Range("H71").Value = 6138.8
$ws.Range("I71").Value = 6173.5
$ws.Range("M71").Value = -27123.5
$ws.Range("K71").Value = 30867.5
$ws.Range("I82").Value = 6916.375
$ws.Range("M82").Value = -6555.375
$ws.Range("K82").Value = 6916.375
$ws.Range("H82").Value = 6729.625
$ws.Range("M85").Value = -5668.375
$ws.Range("K85").Value = 6916.375
$ws.Range("H85").Value = 6729.625
$ws.Range("I85").Value = 6916.375
$ws.Range("H136").Value = 13340572
$ws.Range("I136").Value = 5078.0586
$ws.Range("M136").Value = -12684.1758
$ws.Range("K136").Value = 15234.1758
$ws.Range("N140").Value = -60360
$ws.Range("L140").Value = 50000
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L41").Value = 11957.333
$ws.Range("H41").Value = 12392
$ws.Range("J41").Value = 11957.333
$ws.Range("N41").Value = -12737.333
$ws.Range("I49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("J49").Value = 54000
$ws.Range("N49").Value = -54460
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 54000
$ws.Range("H49").Value = 54000
$ws.Range("L93").Value = 199950
$ws.Range("H93").Value = 199950
$ws.Range("J93").Value = 199950
$ws.Range("N93").Value = -204942
$ws.Range("L99").Value = 149278.25
$ws.Range("H99").Value = 149278.25
$ws.Range("J99").Value = 149278.25
$ws.Range("N99").Value = -155268.25
$ws.Range("H132").Value = 4359.1577
$ws.Range("I132").Value = 3319.9285
$ws.Range("M132").Value = -7429.7855
$ws.Range("J132").Value = 7269
$ws.Range("N132").Value = -26867
$ws.Range("K132").Value = 9959.7855
$ws.Range("L132").Value = 21807
$ws.Range("H136").Value = 4959.375
$ws.Range("I136").Value = 4236.0464
$ws.Range("M136").Value = -10158.1392
$ws.Range("K136").Value = 12708.1392
